$wb = $excel.ActiveWorkbook

# --- Sheet "0" (second sheet) ---
$wsSim = $wb.Worksheets.Item("0")
$wsSim.Activate()

# Select the whole of row 11 (mirrors clicking the row header for row 11)
$wsSim.Rows.Item(11).Select()

# --- Sheet "info" (first sheet) ---
$wsInfo = $wb.Worksheets.Item("info")
$wsInfo.Activate()

# Populate row 15: plot / end_time / 2500 / micros
$wsInfo.Range("A15").Value = "plot"
$wsInfo.Range("B15").Value = "end_time"
$wsInfo.Range("C15").Value = 2500
$wsInfo.Range("D15").Value = "micros"

# Move the active selection to A16 (as recorded after the edit), and
# leave "info" as the active/selected sheet
$wsInfo.Range("A16").Select()
